$d = $word.ActiveDocument

# Each entry: literal text to find -> literal text to replace with,
# using ^l (manual line break => <w:br/>) at the desired split point(s).
# MatchWildcards is left $false so characters like "." and "-" are treated
# literally (no escaping needed).

$replacements = @(
    # Paragraph "Objetivos"
    @("específicas.Integrar", "específicas.^lIntegrar"),
    @("respectivo projeto. - Incentivar", "respectivo projeto. ^l- Incentivar"),

    # Paragraph "Programa"
    @("decisões de saída2.", "decisões de saída^l2."),
    @("inovação tecnológica3.", "inovação tecnológica^l3."),
    @("arquiteturas 4.", "arquiteturas ^l4."),
    @("startups 5.", "startups ^l5."),
    @("pitch de negócio6.", "pitch de negócio^l6."),

    # Paragraph "Método:" run (Avaliação section)
    @("dentre outros.Os alunos", "dentre outros.^lOs alunos"),
    @("profissão.Cada grupo", "profissão.^lCada grupo"),
    @("projeto.As aulas", "projeto.^lAs aulas"),

    # Paragraph "Critério:" run (Avaliação section)
    @("dentre outros.O detalhamento", "dentre outros.^lO detalhamento"),

    # Paragraph "Bibliografia"
    @("LTC, 2017.- BROCKMAN", "LTC, 2017.^l- BROCKMAN"),
    @("LTC, 2010.- CAVALCANTI", "LTC, 2010.^l- CAVALCANTI"),
    @("Saraiva, 2016.- FINOCCHIO", "Saraiva, 2016.^l- FINOCCHIO"),
    @("Saraiva, 2020.- CAMARGO", "Saraiva, 2020.^l- CAMARGO"),
    @("Saraiva, 2019.- BRANCO", "Saraiva, 2019.^l- BRANCO"),
    @("Universitária, 2016- OSTERWALDER", "Universitária, 2016^l- OSTERWALDER")
)

foreach ($pair in $replacements) {
    $find = $pair[0]
    $replace = $pair[1]

    $rng = $d.Content
    $ok = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "FAILED to find: $find"
    }
}

Write-Output "done"
